$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Regolatore" table, rows 16-17.
# Cell values are written in this specific order so that the resulting
# shared-string table indices line up with the authored workbook
# (A16, A17, C16, D16, E16, B16, F16, G16, then the numeric/formula row).
$ws.Range("A16").Value = "Regolatore"
$ws.Range("A17").Value = "Vo = Vref*(1+r2/r1)"
$ws.Range("C16").Value = "vref"
$ws.Range("D16").Value = "r2"
$ws.Range("E16").Value = "r1"
$ws.Range("B16").Value = "vo"
$ws.Range("F16").Value = "r2/r1"
$ws.Range("G16").Value = "vo/vref -1"

$ws.Range("B17").Value = 12
$ws.Range("C17").Value = 1.25
$ws.Range("D17").Formula = "=4.7+3.9"
$ws.Range("E17").Value = 1
$ws.Range("F17").Formula = "=D17/E17"
$ws.Range("G17").Formula = "=B17/C17 - 1"

# Leave the new row selected, matching the author's final cursor position.
$null = $ws.Range("D17").Select()
